# Update recomputed "Return_with_prediction" (G), "return_pct_change" (H),
# and the single "mean_return_pct_change" (I2) values on Sheet1, rows 2-29.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02579693719735697
$ws.Range("H2").Value = -46.51336305070044
$ws.Range("I2").Value = 32.67773463963508

$ws.Range("G3").Value = 0.05919989463311003
$ws.Range("H3").Value = 54.34081560165558

$ws.Range("G4").Value = -0.4618787912001062
$ws.Range("H4").Value = -2.012245335705497

$ws.Range("G5").Value = -0.4604397932635533
$ws.Range("H5").Value = 3.841654472462071

$ws.Range("G6").Value = 0.2556127196486447
$ws.Range("H6").Value = 9.408429361704776

$ws.Range("G7").Value = 0.2702007624851748
$ws.Range("H7").Value = 22.49857791810878

$ws.Range("G8").Value = 0.1608903154894138
$ws.Range("H8").Value = -3.5466992457424

$ws.Range("G9").Value = 0.1849491562926063
$ws.Range("H9").Value = 7.52149672714513

$ws.Range("G10").Value = -0.01274980334362525
$ws.Range("H10").Value = -169.0126001018052

$ws.Range("G11").Value = 0.006291054214953821
$ws.Range("H11").Value = 142.9268324263937

$ws.Range("G12").Value = 0.1350098484745858
$ws.Range("H12").Value = -1.2524853796653

$ws.Range("G13").Value = 0.1384185323755514
$ws.Range("H13").Value = 11.05681768708808

$ws.Range("G14").Value = 0.2582908278018418
$ws.Range("H14").Value = 4.430193756878468

$ws.Range("G15").Value = 0.2666271076568411
$ws.Range("H15").Value = 5.520920841948622

$ws.Range("G16").Value = 0.1381523603370262
$ws.Range("H16").Value = -9.98044137211763

$ws.Range("G17").Value = 0.1475216427201645
$ws.Range("H17").Value = -2.313380189848881

$ws.Range("G18").Value = -0.01658167031579169
$ws.Range("H18").Value = -1.26268851634531

$ws.Range("G19").Value = 0.006315618930864462
$ws.Range("H19").Value = 850.0165762546569

$ws.Range("G20").Value = 0.1322789828876194
$ws.Range("H20").Value = -4.603411204480114

$ws.Range("G21").Value = 0.1418334602758132
$ws.Range("H21").Value = -0.8812464031748179

$ws.Range("G22").Value = 0.1679734769652912
$ws.Range("H22").Value = -9.797102949731416

$ws.Range("G23").Value = 0.187742138699867
$ws.Range("H23").Value = 4.614739331437291

$ws.Range("G24").Value = -0.09216080852107984
$ws.Range("H24").Value = 2.369193873819267

$ws.Range("G25").Value = -0.08942848138610564
$ws.Range("H25").Value = 10.21998908295414

$ws.Range("G26").Value = 0.2331482586053447
$ws.Range("H26").Value = 1.30694803132181

$ws.Range("G27").Value = 0.2485247810341548
$ws.Range("H27").Value = 6.859177569310923

$ws.Range("G28").Value = 0.06968555133095868
$ws.Range("H28").Value = 18.50626676108146

$ws.Range("G29").Value = 0.07814666749210229
$ws.Range("H29").Value = 10.71360396113254
